$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("1.523068561977681", "0.03943062122654339", "0.04763032257267064", "6.00853074385131", "0.002668494548193342", "0.5464607382105982", "1.106373297567018", "0.1357192325528929")
    3 = @("1.516457162708207", "0.03434834991520574", "0.04845170421691236", "5.812542000581487", "0.002673614036068076", "0.522366911553604", "1.089308238722367", "0.1414356777043659")
    4 = @("1.51367713576272", "0.03122369891895005", "0.04901043453154408", "5.692770281592175", "0.002676920212048881", "0.5074412865761389", "1.07977379261672", "0.145188388202051")
    5 = @("1.512865755162636", "0.0299492340242864", "0.04925176388660546", "5.644100068471403", "0.002678308585941767", "0.5013249584952035", "1.076125270785639", "0.1467785608128835")
    6 = @("1.512750438171935", "0.02973753781638777", "0.04929265932978666", "5.636026601210119", "0.00267854160982547", "0.5003072679002543", "1.075533729087255", "0.1470462827100674")
    7 = @("1.513664891731537", "0.0312065158154553", "0.04901363399569725", "5.692113346524224", "0.002676938769674091", "0.5073589384585091", "1.079723628971323", "0.1452095873789343")
    8 = @("1.520523270846695", "0.03767904090902618", "0.0479022222754768", "5.940833458969564", "0.002670226048301696", "0.5381800526206177", "1.100293157150531", "0.1376397880543232")
    9 = @("1.544138340703881", "0.05034451374599769", "0.04615606416521878", "6.433299520019659", "0.002658347433594435", "0.5976137506560377", "1.148140935926193", "0.1247280528290959")
    10 = @("1.567712499916183", "0.05964181614049835", "0.04513970525804289", "6.798357308289013", "0.002650394255884064", "0.6407252275300834", "1.18791435773133", "0.1164284442196069")
    11 = @("1.579795068222296", "0.06387156007068029", "0.04473574130287261", "6.965216740502854", "0.002646942228610179", "0.6602303345382268", "1.207020954458216", "0.1129122829635243")
    12 = @("1.584566237504049", "0.06547344220777518", "0.04459121027112722", "7.028521975997876", "0.002645658740813959", "0.6676020780883221", "1.214402535931413", "0.1116182695393633")
    13 = @("1.583529967574179", "0.06512843876315344", "0.04462196144358721", "7.014882676448053", "0.002645934110094919", "0.6660150680158949", "1.212806262661786", "0.1118952885690883")
    14 = @("1.580183669426219", "0.06400334378572836", "0.04472368120472403", "6.970422481608352", "0.002646836160602906", "0.6608370960737489", "1.207625305513659", "0.1128050714725806")
    15 = @("1.578159471620694", "0.06331421630819989", "0.04478708818546018", "6.943205005756397", "0.002647391777651881", "0.6576635880552431", "1.204470893808121", "0.1133672263404133")
    16 = @("1.566950247410688", "0.05936540969408099", "0.04516728397742575", "6.787469046771889", "0.00265062318081191", "0.6394484562576395", "1.186686140460353", "0.1166634684279382")
    17 = @("1.560421992337155", "0.05694311968925092", "0.04541550750283641", "6.692136802562402", "0.002652647938068957", "0.6282474188943468", "1.176035759036267", "0.1187521667535014")
    18 = @("1.556794961868633", "0.05554990870557219", "0.04556377158514557", "6.63737852570506", "0.002653828149519515", "0.6217947980229326", "1.170005297079825", "0.119977932968431")
    19 = @("1.555588857953921", "0.05507819229852373", "0.04561491320409239", "6.618850919921243", "0.002654230436954906", "0.6196082928358777", "1.167979845978863", "0.1203971419242764")
    20 = @("1.561103701517226", "0.05720097302832983", "0.04538851490515228", "6.702277354779596", "0.002652430783009141", "0.629440825200362", "1.177159637155484", "0.1185272944122415")
    21 = @("1.581161242039656", "0.0643338058866334", "0.04469357417919362", "6.983478242073261", "0.002646570563928651", "0.6623583760164422", "1.209143102095993", "0.1125368274912866")
    22 = @("1.595411262726827", "0.06899657111361535", "0.0442886147557946", "7.167956470064041", "0.002642878767498148", "0.6837882871829777", "1.230899345951116", "0.1088402639851713")
    23 = @("1.587701194469133", "0.06650783010977079", "0.04450022979330193", "7.069431441068446", "0.002644836549330525", "0.6723580913541127", "1.219209360236448", "0.1107931303957761")
    24 = @("1.560795107878278", "0.0570843994261736", "0.04540070095233162", "6.697692657312359", "0.002652528908482615", "0.6289013267711425", "1.176651243512623", "0.1186288815064831")
    25 = @("1.536658821666634", "0.04692019653680291", "0.04658181274767337", "6.299531282200405", "0.002661424311464301", "0.5816370202915238", "1.134388595988071", "0.1280132563442749")
}

$cols = @("B","C","E","F","G","J","K","L")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = [double]$vals[$i]
    }
}

$wb.Save()
